# Add Join Game requirements to the requirements tracking sheet.
$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# New rows of data: Row#, A (id), B (group), C (requirement text), D (implemented-in text)
$rows = @(
    @(53, 50, 26, "The program shall display a waitlist of available players ready to play with their usernames and IDs in the order that the players joined.", "joingame/Waitlist.js/displayWaitlist()"),
    @(54, 51, 26, "The program shall display a button for a game against another user.", "webapp/index.html"),
    @(55, 52, 26, "The program shall display a button for a game against a bot.", "webapp/index.html"),
    @(56, 53, 26, "The program shall update the waitlist as players login and join/leave a game.", "joingame/Waitlist.js/add(), joingame/Waitlist.js/remove()"),
    @(57, 54, 26, "The program shall display a notification with updates to the waitlist.", "joingame/DisplayNotification.js/displayNotification()"),
    @(58, 55, 26, "The program will display a notification to the user when they are next to join a game.", "joingame/DisplayNotification.js/displayNotification()"),
    @(59, 56, 26, "The program shall have display a button to spectate a game of 2 bots playing against each other. ", "joingame/MatchMaker.js/requestSpectateBotVsBot()"),
    @(60, 57, 26, "The program will display players that choose to replay and new players at the bottom of the waitlist.", "joingame/Waitlist.js/displayWaitlist()"),
    @(61, 58, 26, "The program will load usernames and IDs in waitlist and display them within 2 seconds of logging in/joining a game.", "joingame/Waitlist.js/displayWaitlist()"),
    @(62, 59, 26, "The web interface of the Join Game component will be written in HTML5.", "webapp/index.html"),
    @(63, 60, 26, "The program will be compatible with the web browsers Chrome, Firefox, and Safari.", "joingame"),
    @(64, 61, 26, "The interactive and dynamic functionality within the component will be written in JavaScript.", "joingame/Communication.js, Data.js, Displaynotification.js, MatchMaking.js, Waitlist.js"),
    @(65, 62, 26, "When the button to play against another player is pressed, the program will send the player information along with the game mode to the Page Manager component to initalize the pairing up process to begin a game.", "joingame/MatchMaker.js/requestPlayerMatch()"),
    @(66, 63, 26, "When the button to play against a bot is pressed, the program will send the player information along with the game mode to the Page Manager component to intialize the process to start a game against a bot.", "joingame/MatchMaker.js/requestBotMatch()"),
    @(67, 64, 26, "When a player(s) requests to join a game, the program removes their name from the waitlist.", "joingame/Waitlist.js/remove()"),
    @(68, 65, 26, "When a player(s) logins successfully, the program adds their name to the waitlist.", "joingame/Waitlist.js/add()"),
    @(69, 66, 26, "The program will receive player information (username and ID) from the Page Manager through a websocket.", "joingame/Communication.js/handlePlayerData()"),
    @(70, 67, 26, "The program will create a player instance based on the provided data from the Page Manager.", "joingame/Data.js/constructor()"),
    @(71, 68, 26, "The program will send player information (username and ID) and their selected game mode to the Page Manager through a websocket.", "joingame/Communication.js/sendPlayerAttributes()"),
    @(72, 69, 26, "When the button to spectate a game of two bots playing against each other, the program will send players' username and ID along with the game mode to Page Manager to intialize the process to start the game", "joingame/MatchMaker.js/requestSpectateBotVsBot()")
)

foreach ($row in $rows) {
    $r = $row[0]
    $ws.Cells.Item($r, 1).Value = $row[1]
    $ws.Cells.Item($r, 2).Value = $row[2]
    $ws.Cells.Item($r, 3).Value = $row[3]
    $ws.Cells.Item($r, 4).Value = $row[4]
}

# Update sheet view to match the post-edit scroll position/zoom.
$ws.Application.ActiveWindow.Zoom = 75
